# Append 4 new daily rows (189-192) to Sheet1, covering 2020-11-04 .. 2020-11-07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(44139, 904739, 138298, 766132, 309, 2239, 131768, 239, 27, 51, 0, 2, 846, 1189, 39, 123, 42),
    @(44140, 909225, 138567, 770304, 354, 2241, 131987, 235, 25, 63, 0, 2, 847, 1190, 39, 123, 42),
    @(44141, 912762, 138768, 773785, 209, 2242, 132059, 233, 34, 48, 0, 0, 847, 1191, 39, 123, 42),
    @(44142, 916792, 139011, 777548, 233, 2243, 132403, 222, 25, 48, 0, 0, 847, 1192, 39, 123, 42)
)

$startRow = 189
$templateRow = 165   # existing row already formatted like the new rows (style 10/11)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    # Reuse the existing cell styles (date format / right alignment) instead
    # of fabricating brand-new style entries.
    $ws.Range("A$templateRow`:Q$templateRow").Copy() | Out-Null
    $ws.Range("A$row`:Q$row").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false
